# New crime data collected — weekly refresh of the 63rd Precinct CompStat
# report (week-of bump + updated Week-to-Date / 28-Day / YTD / 2-Year crime
# figures).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Header: volume/number and reporting-week dates ----
$ws.Range("A8").Value = "Volume 30   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/16/2023  Through  10/22/2023"

# ---- Numeric value updates (counts / % changes) ----
$ws.Range("M15").Value = -37.5
$ws.Range("N15").Value = -47.368421052631
$ws.Range("F16").Value = 4
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = -42.857142857142
$ws.Range("I16").Value = 92
$ws.Range("K16").Value = -4.166666666666
$ws.Range("L16").Value = 43.75
$ws.Range("M16").Value = -48.603351955307
$ws.Range("N16").Value = -83.512544802867
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 200
$ws.Range("F17").Value = 10
$ws.Range("G17").Value = 11
$ws.Range("H17").Value = -9.090909090909
$ws.Range("I17").Value = 136
$ws.Range("J17").Value = 141
$ws.Range("K17").Value = -3.54609929078
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = 32.038834951456
$ws.Range("N17").Value = -47.892720306513
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 9
$ws.Range("H18").Value = -22.222222222222
$ws.Range("I18").Value = 77
$ws.Range("J18").Value = 93
$ws.Range("K18").Value = -17.204301075268
$ws.Range("L18").Value = -11.494252873563
$ws.Range("M18").Value = -66.810344827586
$ws.Range("N18").Value = -91.773504273504
$ws.Range("C19").Value = 11
$ws.Range("D19").Value = 12
$ws.Range("E19").Value = -8.333333333333
$ws.Range("F19").Value = 45
$ws.Range("G19").Value = 48
$ws.Range("H19").Value = -6.25
$ws.Range("I19").Value = 531
$ws.Range("J19").Value = 421
$ws.Range("K19").Value = 26.128266033254
$ws.Range("L19").Value = 76.41196013289
$ws.Range("M19").Value = 32.089552238806
$ws.Range("N19").Value = -0.187969924812
$ws.Range("C20").Value = 7
$ws.Range("D20").Value = 1
$ws.Range("E20").Value = 600
$ws.Range("F20").Value = 16
$ws.Range("H20").Value = 14.285714285714
$ws.Range("I20").Value = 106
$ws.Range("J20").Value = 103
$ws.Range("K20").Value = 2.912621359223
$ws.Range("L20").Value = 68.253968253968
$ws.Range("M20").Value = -19.696969696969
$ws.Range("N20").Value = -95.261510952168
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 15
$ws.Range("E21").Value = 66.666666666666
$ws.Range("F21").Value = 82
$ws.Range("H21").Value = -8.888888888888
$ws.Range("I21").Value = 955
$ws.Range("J21").Value = 862
$ws.Range("K21").Value = 10.788863109048
$ws.Range("L21").Value = 43.393393393393
$ws.Range("M21").Value = -10.328638497652
$ws.Range("N21").Value = -79.043230195303
$ws.Range("I23").Value = 24
$ws.Range("K23").Value = -11.111111111111
$ws.Range("L23").Value = 4.347826086956
$ws.Range("M23").Value = -4
$ws.Range("C24").Value = 27
$ws.Range("D24").Value = 24
$ws.Range("E24").Value = 12.5
$ws.Range("G24").Value = 94
$ws.Range("H24").Value = 29.787234042553
$ws.Range("I24").Value = 1052
$ws.Range("J24").Value = 924
$ws.Range("K24").Value = 13.852813852813
$ws.Range("L24").Value = 59.635811836115
$ws.Range("M24").Value = 29.556650246305
$ws.Range("C25").Value = 4
$ws.Range("D25").Value = 7
$ws.Range("E25").Value = -42.857142857142
$ws.Range("F25").Value = 15
$ws.Range("G25").Value = 20
$ws.Range("H25").Value = -25
$ws.Range("I25").Value = 216
$ws.Range("J25").Value = 189
$ws.Range("K25").Value = 14.285714285714
$ws.Range("L25").Value = 6.403940886699
$ws.Range("M25").Value = -24.210526315789
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 4
$ws.Range("H27").Value = -50
$ws.Range("I27").Value = 24
$ws.Range("J27").Value = 27
$ws.Range("K27").Value = -11.111111111111
$ws.Range("L27").Value = 14.285714285714

# ---- Cells that went from a real count to "no data" (text "0"/"***.*") ----
# Use a leading apostrophe so the engine stores these as text (matching the
# original "0" / "***.*" placeholder cells elsewhere in the sheet), then
# re-stamp the format from a neighboring placeholder cell (style 14) so the
# number format matches exactly instead of inheriting the old numeric style.
$textCells = "D14","E14","F15","D16","E16","F26","D28","E28","D29","E29"
$textValues = @{
  "D14" = "'0"
  "E14" = "'***.*"
  "F15" = "'0"
  "D16" = "'0"
  "E16" = "'***.*"
  "F26" = "'0"
  "D28" = "'0"
  "E28" = "'***.*"
  "D29" = "'0"
  "E29" = "'***.*"
}
$ws.Range("C14").Copy() | Out-Null
foreach ($addr in $textCells) {
  $ws.Range($addr).Value = $textValues[$addr]
  $ws.Range($addr).PasteSpecial(-4122) | Out-Null
}

# ---- Cells that went from "no data" (text) to a real numeric count ----
# Re-stamp with the numeric-count style (style 15, e.g. I23) after writing
# the number so the format matches the other numeric cells in the row.
$ws.Range("I23").Copy() | Out-Null
$ws.Range("C23").Value = 3
$ws.Range("C23").PasteSpecial(-4122) | Out-Null
$ws.Range("F23").Value = 3
$ws.Range("F23").PasteSpecial(-4122) | Out-Null

$ws.Range("D27").Copy() | Out-Null
$ws.Range("C27").Value = 1
$ws.Range("C27").PasteSpecial(-4122) | Out-Null
